$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PRINCIPAL value in C2
$ws.Range("C2").Value = 100000

# Apply a number format (comma separated with 2 decimals) to match numFmtId 4
$ws.Range("C2").NumberFormat = "#,##0.00"

# Adjust column C width to match bestFit width (closest achievable value to 10.140625)
$ws.Columns.Item(3).ColumnWidth = 9.307291666666666

# Update the selection to J17
$ws.Range("J17").Select()
